$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (y_1_forecast) and column E (y_1) values per the bugfix
$ws.Range("C2").Value = 0.7771450785698075
$ws.Range("E2").Value = 0.4748521911469794

$ws.Range("C3").Value = 0.9995687521967556
$ws.Range("E3").Value = 0.02446583860156171

$ws.Range("C4").Value = 0.5978820435291077
$ws.Range("E4").Value = 0.7697980859487474

$ws.Range("C5").Value = -0.5061359875450311
$ws.Range("E5").Value = -0.362324052998142

$ws.Range("C6").Value = -0.2706540469742502
$ws.Range("E6").Value = -0.06282556559906727

$ws.Range("C7").Value = 0.006024133679316535
$ws.Range("E7").Value = 0.01247916696665019

$ws.Range("C8").Value = 0.3981709080043139
$ws.Range("E8").Value = 0.1740860482467133

$ws.Range("C9").Value = -0.2945738319855118
$ws.Range("E9").Value = -0.03768624985649449

$ws.Range("C10").Value = -0.03047919532177534
$ws.Range("E10").Value = -0.1249617237519152

$ws.Range("C11").Value = 0.1932702877606385
$ws.Range("E11").Value = -0.250093582508859

$ws.Range("C12").Value = -0.1256133802673975
$ws.Range("E12").Value = -0.02501876407304815

$ws.Range("C13").Value = -0.2736870064301455
$ws.Range("E13").Value = -0.11321783823105

$ws.Range("C14").Value = -0.2932081122163033
$ws.Range("E14").Value = -0.1126446518617819

$ws.Range("C15").Value = -0.03584227163500042
$ws.Range("E15").Value = -0.1719585843969917

$ws.Range("C16").Value = 1.323454226677478
$ws.Range("E16").Value = 0.7858329241748896

$ws.Range("C17").Value = 0.8831924739260089
$ws.Range("E17").Value = 0.8355283619122744

$ws.Range("C18").Value = -1.081515348061801
$ws.Range("E18").Value = -0.07003400812273242

$ws.Range("C19").Value = 0.4626514211933497
$ws.Range("E19").Value = -0.4747835872719319
